$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime" timestamps
# for the first data row (row 2) on the zh-cn and de-de report sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 10:21:38"
$wsZhCn.Range("G2").Value = "2016-01-08 10:22:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 10:21:47"
$wsDeDe.Range("G2").Value = "2016-01-08 10:22:39"
